$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IPAR")

# 1. Insert a brand-new column before column D. Everything that used to live in
#    columns D:K (the 2017..2011 fiscal-year data) shifts right to E:L, making
#    room for a new "2018" column at D.
$ws.Columns("D").Insert()

# 2. The freshly inserted column D is blank and picked up formatting from the
#    column to its left (C). Re-copy number formats from column E (which now
#    holds what used to be column D) onto the new column D so the new 2018
#    column renders like the rest of the data (dates as dates, numbers with
#    thousands separators, etc.)
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3. Populate the new column D with the 2018 fiscal-year figures (and re-stamp
#    the "Period Ending" header cells with the new date).
$newColumnD = @(
    @(7, 43465, $false),
    @(8, 675600, $false),
    @(9, 296900, $false),
    @(10, 378700, $false),
    @(12, "NA", $true),
    @(13, 0, $false),
    @(14, 0, $false),
    @(15, 0, $false),
    @(17, 580800, $false),
    @(18, 94700, $false),
    @(20, 3700, $false),
    @(21, 109500, $false),
    @(22, 2600, $false),
    @(23, 95900, $false),
    @(24, 26100, $false),
    @(25, 0, $false),
    @(26, 69700, $false),
    @(27, 53800, $false),
    @(28, 0, $false),
    @(29, "NA", $true),
    @(30, 0, $false),
    @(31, 0, $false),
    @(32, -3700, $false),
    @(33, 53800, $false),
    @(34, 0, $false),
    @(35, 53800, $false),
    @(38, 43465, $false),
    @(41, 193100, $false),
    @(42, 67900, $false),
    @(43, 139300, $false),
    @(44, 161000, $false),
    @(45, 8100, $false),
    @(46, 569400, $false),
    @(47, 0, $false),
    @(48, 9800, $false),
    @(49, 204300, $false),
    @(50, 0, $false),
    @(51, 0, $false),
    @(52, 15600, $false),
    @(53, 0, $false),
    @(54, 799200, $false),
    @(57, 58300, $false),
    @(58, 23200, $false),
    @(59, 105500, $false),
    @(60, 187000, $false),
    @(61, 22900, $false),
    @(62, 3500, $false),
    @(63, 0, $false),
    @(64, 0, $false),
    @(65, 0, $false),
    @(66, 351600, $false),
    @(68, 0, $false),
    @(69, 0, $false),
    @(70, 0, $false),
    @(71, 0, $false),
    @(72, 448700, $false),
    @(73, 0, $false),
    @(74, 0, $false),
    @(75, 0, $false),
    @(76, 447600, $false),
    @(77, 0, $false),
    @(80, 43465, $false),
    @(81, 53800, $false),
    @(83, 11000, $false),
    @(84, 0, $false),
    @(85, 0, $false),
    @(86, 0, $false),
    @(87, 0, $false),
    @(88, 0, $false),
    @(89, 63000, $false),
    @(91, -4000, $false),
    @(92, 0, $false),
    @(93, 0, $false),
    @(94, -13600, $false),
    @(96, -26300, $false),
    @(97, 0, $false),
    @(98, 0, $false),
    @(99, 0, $false),
    @(100, -55900, $false),
    @(101, -8700, $false),
    @(102, -15200, $false)
)

foreach ($item in $newColumnD) {
    $r = $item[0]
    $v = $item[1]
    $isText = $item[2]
    $cell = $ws.Range("D$r")
    if ($isText) {
        $cell.Value2 = [string]$v
    } else {
        $cell.Value2 = $v
    }
}

# 4. "Effect Of Exchange Rate Changes " (row 91) was restated across the board,
#    not just shifted - overwrite D91:K91 with the corrected figures.
$ws.Range("D91").Value2 = -4000
$ws.Range("E91").Value2 = -3000
$ws.Range("F91").Value2 = -4800
$ws.Range("G91").Value2 = -4200
$ws.Range("H91").Value2 = -3300
$ws.Range("I91").Value2 = -5000
$ws.Range("J91").Value2 = -9500
$ws.Range("K91").Value2 = -14600
